# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing table, matching the header style already used by the other
# header cells (bold, centered/top, thin border) and filling every data
# row (including the duplicated header row 38) with the season record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AB1, style "1")
# onto the three new header cells so they match the rest of the header
# row exactly, then set their text.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill the season record for every remaining row (2-38) — the workbook
# only covers a single season so every row gets the same record.
$ws.Range("AC2:AC38").Value = 63
$ws.Range("AD2:AD38").Value = 99
$ws.Range("AE2:AE38").Value = 0
